$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$cell = $ws.Cells.Item(2, 4)
$cell.NumberFormat = "@"
$cell.Value = "43.026.03"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(2, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.55%  "
$cell.Style = "Normal"

# Row 3
$cell = $ws.Cells.Item(3, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.366.13"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(3, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +2.35%  "
$cell.Style = "Normal"

# Row 4
$cell = $ws.Cells.Item(4, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(4, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.04%  "
$cell.Style = "Normal"

# Row 5
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = "302.52"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(5, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.40%  "
$cell.Style = "Normal"

# Row 6
$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = "96.03"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(6, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.61%  "
$cell.Style = "Normal"

# Row 7
$cell = $ws.Cells.Item(7, 2)
$cell.NumberFormat = "@"
$cell.Value = "XRP"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(7, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(7, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.504"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(7, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.23%  "
$cell.Style = "Normal"

# Row 8
$cell = $ws.Cells.Item(8, 2)
$cell.NumberFormat = "@"
$cell.Value = "USDC"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(8, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(8, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(8, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.06%  "
$cell.Style = "Normal"

# Row 9
$cell = $ws.Cells.Item(9, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.491"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(9, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.19%  "
$cell.Style = "Normal"

# Row 10
$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = "@"
$cell.Value = "34.15"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(10, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.11%  "
$cell.Style = "Normal"

# Row 11
$cell = $ws.Cells.Item(11, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +3.67%  "
$cell.Style = "Normal"

# Row 12
$cell = $ws.Cells.Item(12, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0786"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(12, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.38%  "
$cell.Style = "Normal"

# Row 13
$cell = $ws.Cells.Item(13, 4)
$cell.NumberFormat = "@"
$cell.Value = "18.39"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(13, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -2.99%  "
$cell.Style = "Normal"

# Row 14
$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = "@"
$cell.Value = "6.75"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(14, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.50%  "
$cell.Style = "Normal"

# Row 15
$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.731.26"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(15, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +2.23%  "
$cell.Style = "Normal"

# Row 16
$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.359.55"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(16, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +2.87%  "
$cell.Style = "Normal"

# Row 17
$cell = $ws.Cells.Item(17, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.796"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(17, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.87%  "
$cell.Style = "Normal"

# Row 18
$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = "@"
$cell.Value = "42.977.97"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(18, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.57%  "
$cell.Style = "Normal"

# Row 19
$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = "@"
$cell.Value = "11.99"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(19, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -1.69%  "
$cell.Style = "Normal"

# Row 20
$cell = $ws.Cells.Item(20, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +2.15%  "
$cell.Style = "Normal"

# Row 21
$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0886"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(21, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.33%  "
$cell.Style = "Normal"

# Row 22
$cell = $ws.Cells.Item(22, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.40%  "
$cell.Style = "Normal"

# Row 23
$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = "234.93"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(23, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.04%  "
$cell.Style = "Normal"

# Row 24
$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.18"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(24, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -4.32%  "
$cell.Style = "Normal"

# Row 25
$cell = $ws.Cells.Item(25, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.01%  "
$cell.Style = "Normal"

# Row 26
$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.43"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(26, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.16%  "
$cell.Style = "Normal"

# Row 27
$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = "24.55"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(27, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +1.14%  "
$cell.Style = "Normal"

# Row 28
$cell = $ws.Cells.Item(28, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.28%  "
$cell.Style = "Normal"

# Row 29
$cell = $ws.Cells.Item(29, 4)
$cell.NumberFormat = "@"
$cell.Value = "9.28"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(29, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +1.85%  "
$cell.Style = "Normal"

# Row 30
$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = "@"
$cell.Value = "31.78"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(30, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.98%  "
$cell.Style = "Normal"

# Row 31
$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(31, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.05%  "
$cell.Style = "Normal"

# Row 32
$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = "@"
$cell.Value = "5.04"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(32, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.68%  "
$cell.Style = "Normal"

# Row 33
$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = "@"
$cell.Value = "17.36"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(33, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -1.33%  "
$cell.Style = "Normal"

# Row 34
$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0718"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(34, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +2.88%  "
$cell.Style = "Normal"

# Row 35
$cell = $ws.Cells.Item(35, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +4.22%  "
$cell.Style = "Normal"

# Row 36
$cell = $ws.Cells.Item(36, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +3.80%  "
$cell.Style = "Normal"

# Row 37
$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = "@"
$cell.Value = "4.37"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(37, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -2.03%  "
$cell.Style = "Normal"

# Row 38
$cell = $ws.Cells.Item(38, 2)
$cell.NumberFormat = "@"
$cell.Value = "WEMIXToken"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(38, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.29"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(38, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -1.63%  "
$cell.Style = "Normal"

# Row 39
$cell = $ws.Cells.Item(39, 2)
$cell.NumberFormat = "@"
$cell.Value = "Monero"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(39, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = "@"
$cell.Value = "123.16"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(39, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -26.00%  "
$cell.Style = "Normal"

# Row 40
$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.79"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(40, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +3.03%  "
$cell.Style = "Normal"

# Row 41
$cell = $ws.Cells.Item(41, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.45%  "
$cell.Style = "Normal"

# Row 42
$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = "@"
$cell.Value = "21.60"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(42, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +2.81%  "
$cell.Style = "Normal"

# Row 43
$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.937.86"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(43, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.64%  "
$cell.Style = "Normal"

# Row 44
$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0279"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(44, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.15%  "
$cell.Style = "Normal"

# Row 45
$cell = $ws.Cells.Item(45, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +2.33%  "
$cell.Style = "Normal"

# Row 46
$cell = $ws.Cells.Item(46, 2)
$cell.NumberFormat = "@"
$cell.Value = "NEARProtocol"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(46, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.73"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(46, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.13%  "
$cell.Style = "Normal"

# Row 47
$cell = $ws.Cells.Item(47, 2)
$cell.NumberFormat = "@"
$cell.Value = "FraxShare"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(47, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = "@"
$cell.Value = "9.19"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(47, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -8.80%  "
$cell.Style = "Normal"

# Row 48
$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.589.51"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(48, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +1.91%  "
$cell.Style = "Normal"

# Row 49
$cell = $ws.Cells.Item(49, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +2.52%  "
$cell.Style = "Normal"

# Row 50
$cell = $ws.Cells.Item(50, 2)
$cell.NumberFormat = "@"
$cell.Value = "MultiversX"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(50, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = "@"
$cell.Value = "51.76"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(50, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -3.03%  "
$cell.Style = "Normal"

# Row 51
$cell = $ws.Cells.Item(51, 2)
$cell.NumberFormat = "@"
$cell.Value = "TrustWalletToken"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(51, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.14"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(51, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +1.52%  "
$cell.Style = "Normal"

Write-Output "Update complete"